$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    3  = @(0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.069511820747223)
    4  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    5  = @(0.01293466051926884, 0.306821227259698, 0.1494219747398047, 10.19245300693656, 10.66163086945533)
    6  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    7  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    8  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    9  = @(0.04271373187048222, 0.306821227259698, 0.7527432677738641, 10.19245300693656, 11.2947312338406)
    10 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 10.19245300693656, 15.28448560880142)
    11 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 18.67282528286833)
    12 = @(0.1190320826869504, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 4.457851494814137)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
